# 24.06.19 Today Sales Updated
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Raju Ahamed")

# Update the date header (both occurrences share the same shared string)
$ws.Range("A4").Value = "Date: 24.06.19"
$ws.Range("A31").Value = "Date: 24.06.19"

# First table (rows 6-11)
$ws.Range("E6").Value = 27
$ws.Range("E7").Value = 86
$ws.Range("E8").Value = 100
$ws.Range("E9").Value = $null
$ws.Range("E10").Value = $null
$ws.Range("E11").Value = $null

# Second table (rows 33-38)
$ws.Range("E33").Value = 27
$ws.Range("E34").Value = 86
$ws.Range("E35").Value = 100
$ws.Range("E36").Value = $null
$ws.Range("E37").Value = $null
$ws.Range("E38").Value = $null

$excel.Calculate()
